$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.786.74'
$ws.Range("E2").Value = '  +7.15%  '
$ws.Range("D3").Value = '3.854.77'
$ws.Range("E3").Value = '  +12.62%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").Value = '''426.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +11.21%  '
$ws.Range("D6").Value = '''131.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +9.01%  '
$ws.Range("D7").Value = '3.848.00'
$ws.Range("E7").Value = '  +12.57%  '
$ws.Range("D8").Value = '''0.611'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.10%  '
$ws.Range("E9").Value = '  -0.12%  '
$ws.Range("D10").Value = '''0.726'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.55%  '
$ws.Range("E11").Value = '  +12.08%  '
$ws.Range("D12").Value = '''0.0000343'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +15.75%  '
$ws.Range("D13").Value = '''41.10'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +8.18%  '
$ws.Range("D14").Value = '''10.29'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +14.57%  '
$ws.Range("D15").Value = '4.476.71'
$ws.Range("E15").Value = '  +13.75%  '
$ws.Range("D16").Value = '''16.01'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +31.85%  '
$ws.Range("D17").Value = '3.884.28'
$ws.Range("E17").Value = '  +12.27%  '
$ws.Range("D18").Value = '''0.138'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.49%  '
$ws.Range("D19").Value = '''19.99'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +9.91%  '
$ws.Range("D20").Value = '67.029.20'
$ws.Range("E20").Value = '  +7.70%  '
$ws.Range("E21").Value = '  +9.17%  '
$ws.Range("D22").Value = '''413.32'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.41%  '
$ws.Range("D23").Value = '''14.98'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +10.69%  '
$ws.Range("D24").Value = '''84.42'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.01%  '
$ws.Range("D25").Value = '''3.05'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +10.12%  '
$ws.Range("D26").Value = '''37.76'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +17.00%  '
$ws.Range("D27").Value = '''9.98'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +15.88%  '
$ws.Range("D28").Value = '''3.26'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +12.54%  '
$ws.Range("D29").Value = '''5.35'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.10%  '
$ws.Range("D30").Value = '''9.06'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +37.39%  '
$ws.Range("D31").Value = '''720.10'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +11.45%  '
$ws.Range("D32").Value = '''13.58'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +17.65%  '
$ws.Range("D33").Value = '''0.123'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +14.71%  '
$ws.Range("E34").Value = '  +6.97%  '
$ws.Range("E35").Value = '  -0.13%  '
$ws.Range("D36").Value = '''39.31'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.29%  '
$ws.Range("E37").Value = '  +3.22%  '
$ws.Range("D38").Value = '''55.88'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.96%  '
$ws.Range("B39").Value = 'NEARProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D39").Value = '''5.34'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +36.15%  '
$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '0.0₃0757'
$ws.Range("E40").Value = '  +25.03%  '
$ws.Range("D41").Value = '''0.0462'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.90%  '
$ws.Range("D42").Value = '''2.89'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +12.45%  '
$ws.Range("E43").Value = '  +1.12%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").Value = '''0.136'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.35%  '
$ws.Range("B45").Value = 'LidoDAOToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D45").Value = '''3.37'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +13.10%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = '''3.20'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.90%  '
$ws.Range("D47").Value = '''0.313'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +15.62%  '
$ws.Range("D48").Value = '''142.26'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.93%  '
$ws.Range("D49").Value = '''2.05'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.42%  '
$ws.Range("D50").Value = '''2.59'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.36%  '
$ws.Range("D51").Value = '''2.82'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +9.13%  '
